# Pedido 69249e010f48bb0af0483b7d
# - Row 18 (Ignacio Rodriguez): drop the leftover blank placeholder cells
#   for Optimizador/Unidades Optimizador (F:G) and Baterias/Unidades
#   Baterias/Cargador VE (J:L) - the row has no optimizer, battery or EV
#   charger, so those cells go back to being unset.
# - Add a brand-new row 19 for Marcos Cortecero Torres's order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: clear the now-empty columns ---------------------------------
$ws.Range("F18:G18").ClearContents()
$ws.Range("J18:L18").ClearContents()

# --- Row 19: new product/order line --------------------------------------
$ws.Range("A19").Value = 2804
$ws.Range("B19").Value = "Marcos Cortecero Torres"
$ws.Range("C19").Value = "Estructura coplanar NOVOTEGRA"
$ws.Range("D19").Value = "MODULO FV JA SOLAR 535WP BLACK FRAME BIFACIAL 120 CELDAS"

# Columns E:L mirror the sheet's convention of storing these as text
# (e.g. "6", "1") rather than numbers, so force text formatting before
# writing the values.
$ws.Range("E19:L19").NumberFormat = "@"

$ws.Range("E19").Value = "14"
$ws.Range("F19").Value = ""
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = "GOODWE GW6000-ES-20 híbrido monofásico"
$ws.Range("I19").Value = "1"
$ws.Range("J19").Value = "GOODWE Batería Lynx Home U G3 5,12 kWh"
$ws.Range("K19").Value = "1"
$ws.Range("L19").Value = ""

# Drop the temporary "Text" number format again so the cells are left
# without an explicit style override, matching the rest of the sheet.
$ws.Range("E19:L19").Style = "Normal"

$ws.Range("M19").Value = "Sí"
$ws.Range("N19").Value = "2025-09-25T07:50:43.054Z"
